# Updates for Kennedy usage:
# On slide 28 ("unicycler_prerun" demo), the scratch-directory path no
# longer includes the "bioinf/${USER}" segment - it is now just
# "~/scratch/<workshop>/...". Remove that segment from the first
# paragraph of the content placeholder.

$p = $ppt.ActivePresentation

# Locate the slide and the shape that holds the "cd ~/scratch/..." command.
$targetSlideIndex = 28
$s = $p.Slides.Item($targetSlideIndex)

$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTextFrame) {
        if ($candidate.TextFrame.TextRange.Text -like "*bioinf*") {
            $shape = $candidate
            break
        }
    }
}

if ($shape -eq $null) {
    throw "Could not find the shape containing the scratch path text on slide $targetSlideIndex"
}

$tr = $shape.TextFrame.TextRange
$para1 = $tr.Paragraphs(1)

# Build the substring to remove without relying on literal "$" interpolation:
# "bioinf" + "/${USER} /"
$dollar = [char]36
$marker = "bioinf" + "/" + $dollar + "{USER} /"

$fullText = $para1.Text
$markerStart0 = $fullText.IndexOf($marker)
if ($markerStart0 -lt 0) {
    throw "Expected marker text not found in paragraph: $marker"
}

# Characters() is 1-based.
$delRange = $para1.Characters($markerStart0 + 1, $marker.Length)
$delRange.Text = ""

# At this point the paragraph reads "cd ~/scratch/genome_assembly_workshop/unicycler_prerun"
# held as a single leading run "cd ~/scratch/". Split it into "cd " and
# "~/scratch/" (matching the target authoring) by re-touching just the
# "~/scratch/" portion so the engine emits it as its own run.
$scratchMarker = "~/scratch/"
$afterDeleteText = $para1.Text
$scratchStart0 = $afterDeleteText.IndexOf($scratchMarker)
if ($scratchStart0 -lt 0) {
    throw "Expected '~/scratch/' text not found after deletion"
}

$scratchRange = $para1.Characters($scratchStart0 + 1, $scratchMarker.Length)
$scratchRange.Text = $scratchMarker
